$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 0.8337083923857449
$ws.Range("C2").Value = 0.2513404471783076
$ws.Range("D2").Value = 0.009775792956261853
$ws.Range("E2").Value = 0.4254114885864624
$ws.Range("F2").Value = 0.4674185242574538
$ws.Range("I2").Value = 0.2905544193914125
$ws.Range("N2").Value = 0.7993426859598713
$ws.Range("O2").Value = 1.464486666278532

# Row 3
$ws.Range("B3").Value = 0.7296308337800497
$ws.Range("C3").Value = 0.2205598746211024
$ws.Range("D3").Value = 0.008679553284519415
$ws.Range("E3").Value = 0.3710218854530609
$ws.Range("F3").Value = 0.4560650490907463
$ws.Range("I3").Value = 0.2932413434292371
$ws.Range("N3").Value = 0.8024938452141441
$ws.Range("O3").Value = 1.440259507361134

# Row 4
$ws.Range("B4").Value = 0.6655983624438591
$ws.Range("C4").Value = 0.2015752291992499
$ws.Range("D4").Value = 0.008002938075605925
$ws.Range("E4").Value = 0.3377214407544784
$ws.Range("F4").Value = 0.4495041414594994
$ws.Range("I4").Value = 0.295163531058261
$ws.Range("N4").Value = 0.8047950743569103
$ws.Range("O4").Value = 1.4267284468892

# Row 5
$ws.Range("B5").Value = 0.6394734554434365
$ws.Range("C5").Value = 0.1938175921240486
$ws.Range("D5").Value = 0.007726350409907212
$ws.Range("E5").Value = 0.3241732045882486
$ws.Range("F5").Value = 0.4469332072758974
$ws.Range("I5").Value = 0.2960150615022421
$ws.Range("N5").Value = 0.8058251047535592
$ws.Range("O5").Value = 1.421550883858998

# Row 6
$ws.Range("B6").Value = 0.6351335866031036
$ws.Range("C6").Value = 0.1925281654253581
$ws.Range("D6").Value = 0.007680371904793759
$ws.Range("E6").Value = 0.3219248010470466
$ws.Range("F6").Value = 0.4465124945792738
$ws.Range("I6").Value = 0.2961605699614296
$ws.Range("N6").Value = 0.806001716652986
$ws.Range("O6").Value = 1.420711426733703

# Row 7
$ws.Range("B7").Value = 0.6652461574408051
$ws.Range("C7").Value = 0.2014706927276677
$ws.Range("D7").Value = 0.007999211374542625
$ws.Range("E7").Value = 0.3375386380836005
$ws.Range("F7").Value = 0.4494690537591381
$ws.Range("I7").Value = 0.2951747392282904
$ws.Range("N7").Value = 0.8048085919793095
$ws.Range("O7").Value = 1.426657260347412

# Row 8
$ws.Range("B8").Value = 0.7978497558006552
$ws.Range("C8").Value = 0.2407451076477685
$ws.Range("D8").Value = 0.009398553761698736
$ws.Range("E8").Value = 0.40663706359922
$ws.Range("F8").Value = 0.4634184375461388
$ws.Range("I8").Value = 0.2914241671066407
$ws.Range("N8").Value = 0.8003532306202743
$ws.Range("O8").Value = 1.455853091226174

# Row 9
$ws.Range("B9").Value = 1.056828702488701
$ws.Range("C9").Value = 0.3170807680723442
$ws.Range("D9").Value = 0.01211383117156828
$ws.Range("E9").Value = 0.5429856130370609
$ws.Range("F9").Value = 0.4940502139279914
$ws.Range("I9").Value = 0.2862432947631497
$ws.Range("N9").Value = 0.7945188496768907
$ws.Range("O9").Value = 1.52385325336698

# Row 10
$ws.Range("B10").Value = 1.246425377449555
$ws.Range("C10").Value = 0.3727480738894542
$ws.Range("D10").Value = 0.01409015386830248
$ws.Range("E10").Value = 0.6438151142503443
$ws.Range("F10").Value = 0.5185874856690162
$ws.Range("I10").Value = 0.2837798657615132
$ws.Range("N10").Value = 0.7919967240460863
$ws.Range("O10").Value = 1.580482046595961

# Row 11
$ws.Range("B11").Value = 1.332525869809217
$ws.Range("C11").Value = 0.3979822690664605
$ws.Range("D11").Value = 0.01498499671310327
$ws.Range("E11").Value = 0.6898576372912117
$ws.Range("F11").Value = 0.5301989715296287
$ws.Range("I11").Value = 0.2829546248018744
$ws.Range("N11").Value = 0.7912315817689262
$ws.Range("O11").Value = 1.607717883833431

# Row 12
$ws.Range("B12").Value = 1.365107740544488
$ws.Range("C12").Value = 0.4075248725972642
$ws.Range("D12").Value = 0.01532322634214722
$ws.Range("E12").Value = 0.7073201664084365
$ws.Range("F12").Value = 0.5346611084532924
$ws.Range("I12").Value = 0.2826849136235907
$ws.Range("N12").Value = 0.7909967253710448
$ws.Range("O12").Value = 1.618245450402526

# Row 13
$ws.Range("B13").Value = 1.358091672766079
$ws.Range("C13").Value = 0.4054702855544861
$ws.Range("D13").Value = 0.0152504108554794
$ws.Range("E13").Value = 0.7035580537039294
$ws.Range("F13").Value = 0.5336972032175282
$ws.Range("I13").Value = 0.2827410926642244
$ws.Range("N13").Value = 0.7910448658780354
$ws.Range("O13").Value = 1.615968606099898

# Row 14
$ws.Range("B14").Value = 1.335206858251922
$ws.Range("C14").Value = 0.3987676073635384
$ws.Range("D14").Value = 0.01501283579332835
$ws.Range("E14").Value = 0.6912937322958186
$ws.Range("F14").Value = 0.5305647662509045
$ws.Range("I14").Value = 0.282931576132345
$ws.Range("N14").Value = 0.7912111603729244
$ws.Range("O14").Value = 1.608579694783771

# Row 15
$ws.Range("B15").Value = 1.32118627963456
$ws.Range("C15").Value = 0.3946603189524467
$ws.Range("D15").Value = 0.01486723170828697
$ws.Range("E15").Value = 0.6837850897021269
$ws.Range("F15").Value = 0.5286545535408891
$ws.Range("I15").Value = 0.2830538345589986
$ws.Range("N15").Value = 0.7913201663712641
$ws.Range("O15").Value = 1.604081689686353

# Row 16
$ws.Range("B16").Value = 1.240795423983286
$ws.Range("C16").Value = 0.3710971472556253
$ws.Range("D16").Value = 0.01403158724234999
$ws.Range("E16").Value = 0.6408098216017208
$ws.Range("F16").Value = 0.5178377328701629
$ws.Range("I16").Value = 0.2838397678915605
$ws.Range("N16").Value = 0.7920544121841147
$ws.Range("O16").Value = 1.578731959610508

# Row 17
$ws.Range("B17").Value = 1.191439452137672
$ws.Range("C17").Value = 0.356618915856302
$ws.Range("D17").Value = 0.01351785414180284
$ws.Range("E17").Value = 0.614492117360868
$ws.Range("F17").Value = 0.5113174117756216
$ws.Range("I17").Value = 0.2843977906168895
$ws.Range("N17").Value = 0.7926026821403127
$ws.Range("O17").Value = 1.563559830715832

# Row 18
$ws.Range("B18").Value = 1.163037346760632
$ws.Range("C18").Value = 0.3482830426630414
$ws.Range("D18").Value = 0.01322197492827826
$ws.Range("E18").Value = 0.599371184180427
$ws.Range("F18").Value = 0.5076093610925767
$ws.Range("I18").Value = 0.2847465347676881
$ws.Range("N18").Value = 0.7929540094432923
$ws.Range("O18").Value = 1.554971958593853

# Row 19
$ws.Range("B19").Value = 1.153418539989332
$ws.Range("C19").Value = 0.3454592271489503
$ws.Range("D19").Value = 0.01312172847342907
$ws.Range("E19").Value = 0.5942542341585124
$ws.Range("F19").Value = 0.5063611227276823
$ws.Range("I19").Value = 0.2848693759781113
$ws.Range("N19").Value = 0.7930791441920206
$ws.Range("O19").Value = 1.552088025037477

# Row 20
$ws.Range("B20").Value = 1.196694922610163
$ws.Range("C20").Value = 0.358161016943626
$ws.Range("D20").Value = 0.0135725828136799
$ws.Range("E20").Value = 0.6172919782788426
$ws.Range("F20").Value = 0.5120071331272982
$ws.Range("I20").Value = 0.2843355103097842
$ws.Range("N20").Value = 0.7925405949250859
$ws.Range("O20").Value = 1.565160557957284

# Row 21
$ws.Range("B21").Value = 1.341929303930385
$ws.Range("C21").Value = 0.4007367020259949
$ws.Range("D21").Value = 0.01508263458398318
$ws.Range("E21").Value = 0.6948953050097231
$ws.Range("F21").Value = 0.5314830675911963
$ws.Range("I21").Value = 0.2828744627164674
$ws.Range("N21").Value = 0.791160826584516
$ws.Range("O21").Value = 1.610744176116469

# Row 22
$ws.Range("B22").Value = 1.436716553556153
$ws.Range("C22").Value = 0.4284862748506271
$ws.Range("D22").Value = 0.0160658694104896
$ws.Range("E22").Value = 0.745773246523612
$ws.Range("F22").Value = 0.54459151569732
$ws.Range("I22").Value = 0.2821691408395779
$ws.Range("N22").Value = 0.7905789595269823
$ws.Range("O22").Value = 1.641783600532932

# Row 23
$ws.Range("B23").Value = 1.386139281823887
$ws.Range("C23").Value = 0.4136828405773372
$ws.Range("D23").Value = 0.01554144244084199
$ws.Range("E23").Value = 0.718603418774137
$ws.Range("F23").Value = 0.5375603722156512
$ws.Range("I23").Value = 0.2825226460109533
$ws.Range("N23").Value = 0.7908602647480905
$ws.Range("O23").Value = 1.625102487608984

# Row 24
$ws.Range("B24").Value = 1.194319007484694
$ws.Range("C24").Value = 0.3574638708618636
$ws.Range("D24").Value = 0.01354784162219858
$ws.Range("E24").Value = 0.6160261318754294
$ws.Range("F24").Value = 0.5116951837697457
$ws.Range("I24").Value = 0.2843635802579279
$ws.Range("N24").Value = 0.7925685520326624
$ws.Range("O24").Value = 1.564436449441331

# Row 25
$ws.Range("B25").Value = 0.9868843801821754
$ws.Range("C25").Value = 0.2965030198498937
$ws.Range("D25").Value = 0.01138246904334039
$ws.Range("E25").Value = 0.505994874012373
$ws.Range("F25").Value = 0.4854091329161747
$ws.Range("I25").Value = 0.2874102983668045
$ws.Range("N25").Value = 0.7957870164205758
$ws.Range("O25").Value = 1.50429474825043
